$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.883.76'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.894.16'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7829'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.83'
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3139'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.80'
$ws.Range("E9").Value = '  +2.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07317'
$ws.Range("E10").Value = '  +5.28%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7743'
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.474'
$ws.Range("E13").Value = '  +4.18%  '
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.88'
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.794.11'
$ws.Range("E15").Value = '  -5.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.212'
$ws.Range("E16").Value = '  +5.82%  '
$ws.Range("D17").Value = '29.845.28'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.02'
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007812'
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.123'
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.124.37'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.429'
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.13'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.72'
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.020'
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.437'
$ws.Range("E30").Value = '  +2.56%  '
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.480'
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05560'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.059'
$ws.Range("E34").Value = '  +1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.241'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7527'
$ws.Range("E36").Value = '  +3.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.685'
$ws.Range("E38").Value = '  +2.02%  '
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.802'
$ws.Range("E40").Value = '  +1.20%  '
$ws.Range("D41").Value = '1.139.63'
$ws.Range("E41").Value = '  +12.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4461'
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '73.96'
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.960'
$ws.Range("E44").Value = '  +2.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8530'
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9996'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.890'
$ws.Range("E47").Value = '  +2.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.35'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.063'
$ws.Range("E49").Value = '  +6.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.521'
$ws.Range("E50").Value = '  +2.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.757'
$ws.Range("E51").Value = '  -0.63%  '
